$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.633.64'
$ws.Range('E2').Value = '  +0.49%  '
$ws.Range('D3').Value = '1.637.10'
$ws.Range('E3').Value = '  -0.53%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = "'212.59"
$ws.Range('E5').Value = '  -0.03%  '
$ws.Range('D6').Value = "'0.523"
$ws.Range('E6').Value = '  -1.21%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('E8').Value = '  -2.44%  '
$ws.Range('E9').Value = '  -0.47%  '
$ws.Range('E10').Value = '  -0.27%  '
$ws.Range('D11').Value = "'0.0891"
$ws.Range('E11').Value = '  -0.10%  '
$ws.Range('D12').Value = '1.868.95'
$ws.Range('E12').Value = '  -0.52%  '
$ws.Range('D13').Value = '1.627.28'
$ws.Range('E13').Value = '  -1.11%  '
$ws.Range('E14').Value = '  -0.51%  '
$ws.Range('D15').Value = "'0.557"
$ws.Range('E15').Value = '  -5.33%  '
$ws.Range('E16').Value = '  -0.10%  '
$ws.Range('D17').Value = '27.636.64'
$ws.Range('E17').Value = '  +0.64%  '
$ws.Range('D18').Value = "'229.06"
$ws.Range('E18').Value = '  -1.02%  '
$ws.Range('D19').Value = "'7.75"
$ws.Range('E19').Value = '  +2.34%  '
$ws.Range('E20').Value = '  -0.33%  '
$ws.Range('E21').Value = '  +0.04%  '
$ws.Range('E22').Value = '  -1.26%  '
$ws.Range('D23').Value = "'10.06"
$ws.Range('E23').Value = '  +3.47%  '
$ws.Range('E24').Value = '  -1.56%  '
$ws.Range('D25').Value = "'150.21"
$ws.Range('E25').Value = '  +1.62%  '
$ws.Range('E26').Value = '  -1.35%  '
$ws.Range('E27').Value = '  -2.19%  '
$ws.Range('E28').Value = '  +0.08%  '
$ws.Range('D29').Value = "'15.58"
$ws.Range('E29').Value = '  -0.61%  '
$ws.Range('E30').Value = '  -0.28%  '
$ws.Range('E31').Value = '  -0.23%  '
$ws.Range('D32').Value = "'3.29"
$ws.Range('E32').Value = '  -0.31%  '
$ws.Range('D33').Value = '1.454.14'
$ws.Range('E33').Value = '  +1.97%  '
$ws.Range('E34').Value = '  -2.56%  '
$ws.Range('E35').Value = '  -1.67%  '
$ws.Range('E36').Value = '  -0.17%  '
$ws.Range('D37').Value = "'0.562"
$ws.Range('E37').Value = '  -1.22%  '
$ws.Range('E38').Value = '  -1.73%  '
$ws.Range('E39').Value = '  -0.15%  '
$ws.Range('D40').Value = "'0.900"
$ws.Range('E40').Value = '  +8.69%  '
$ws.Range('D41').Value = "'69.79"
$ws.Range('E41').Value = '  +7.52%  '
$ws.Range('E42').Value = '  +0.08%  '
$ws.Range('E43').Value = '  -0.65%  '
$ws.Range('D44').Value = "'5.61"
$ws.Range('E44').Value = '  +1.14%  '
$ws.Range('D45').Value = "'2.46"
$ws.Range('E45').Value = '  +0.14%  '
$ws.Range('E46').Value = '  -0.38%  '
$ws.Range('D47').Value = '1.778.70'
$ws.Range('E47').Value = '  -0.55%  '
$ws.Range('E48').Value = '  +1.87%  '
$ws.Range('E49').Value = '  -2.50%  '
$ws.Range('E50').Value = '  +0.76%  '
$ws.Range('D51').Value = "'0.0985"
$ws.Range('E51').Value = '  -1.04%  '
